# Build site at 2023-04-12 14:53:07 UTC
# Applies the LOB1202.xlsx content update:
#  - inserts a new row (the "Docentes responsaveis" professor name row)
#    right after row 12, pushing everything from the old row 13 onward
#    down by one row
#  - rewrites the PT-BR "Objetivos" / "Programa resumido" / "Programa"
#    cells with their real content (they previously held copy/paste
#    leftovers from neighbouring rows)
#  - rewrites "Metodo" / "Criterio" / "Norma de recuperacao" which were
#    similarly shifted by one row in the old sheet
#  - fills in the previously-empty "Bibliografia" cell
#  - narrows column A's width spec so it no longer also covers column B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at 13 (shifts old rows 13-24 down to 14-25) ---
$ws.Rows.Item(13).Insert()

# --- 2. New row 13: professor name, under "Docentes responsaveis:" (A12) ---
$ws.Range("B13").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"
$ws.Range("C13").Value = "7926291 - Célia Regina Tomachuk dos Santos Catuogno"

# --- 3. Row 10 "Objetivos:" — real objectives text (PT) ---
$ws.Range("B10").Value = "Definição de projeto, ciclo de vida de projetos, processos de gerenciamento de projetos, sustentabilidade, estudo de caso"
$ws.Range("C10").Value = "Definição de projeto, ciclo de vida de projetos, processos de gerenciamento de projetos, sustentabilidade, estudo de caso"

# --- 4. Row 14 "Programa resumido:" — same short syllabus text (PT) ---
$ws.Range("B14").Value = "Definição de projeto, ciclo de vida de projetos, processos de gerenciamento de projetos, sustentabilidade, estudo de caso"
$ws.Range("C14").Value = "Definição de projeto, ciclo de vida de projetos, processos de gerenciamento de projetos, sustentabilidade, estudo de caso"

# --- 5. Row 16 "Programa:" — full syllabus text (PT) ---
$programaPt = @"
Definição de projeto e seus principais atributos e características; conceitos do PMBoK (Project managment body of knowledge). Planejamento estratégico. Desenvolvimento Sustentável: O que é desenvolvimento sustentável? Convênios, tratados e políticas de alcance internacional realizado em torno do desenvolvimento sustentável. Os desafios do desenvolvimento sustentável. Processos e metodologia do gerenciamento de projetos ambientais. Ferramentas de planejamento, monitoramento e controle. Estudo dos riscos e problemas comuns na gestão de projetos ambientais. Análise de casos reais envolvendo seleção, administração e desenvolvimento de projetos aplicados à gestão ambiental
"@
$programaPt = $programaPt.TrimEnd("`r","`n")
$ws.Range("B16").Value = $programaPt
$ws.Range("C16").Value = $programaPt

# --- 6. Row 19 "Método:" — evaluation method text ---
$ws.Range("B19").Value = "Avaliação composta por 2 (duas) provas e um projetoPara os alunos que perderem uma das provas será oferecida uma substitutiva no final do semestre, que incluirá toda a matéria da disciplina."
$ws.Range("C19").Value = "Avaliação composta por 2 (duas) provas e um projetoPara os alunos que perderem uma das provas será oferecida uma substitutiva no final do semestre, que incluirá toda a matéria da disciplina."

# --- 7. Row 20 "Critério:" — final grade formula ---
$ws.Range("B20").Value = "Nota final = (nota prova 1 + nota da prova 2 + nota do projeto)/3."
$ws.Range("C20").Value = "Nota final = (nota prova 1 + nota da prova 2 + nota do projeto)/3."

# --- 8. Row 21 "Norma de recuperação:" — make-up exam rule ---
$ws.Range("B21").Value = "Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota recuperação)/2] deverá ser igual ou superior a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com todo o conteúdo da disciplina, sendo que a nota [(nota final do semestre + nota recuperação)/2] deverá ser igual ou superior a 5,0 (cinco)."

# --- 9. Row 22 "Bibliografia:" — bibliography list (was empty) ---
$bibliografia = @"
Bibliografia básica:
VALLE, A.B. (Org.), Fundamentos do gerenciamento de projetos, 3° edição, Editora FGV, 2008
MENEZES, L.C.M., Gestão de projetos, 2° edição, Editora Atlas, 2009
KAHN, M., Gerenciamento de projetos ambientais, E-papers Serviços Editoriais, 2003
VARGAS, R., Manual prático de projeto, 3° edição, Editora Brasport, 2007
Bibliografia complementar:
TORRES, C., Lélis, J.C., Garantia de sucesso em gestão de projetos, Ed. Brasport, 2008
ROCHA, J.S.M., Manual de projetos ambientais, Imprensa Universitária, 1997
Project Management Institute. PMBok, 2013
"@
$bibliografia = $bibliografia.TrimEnd("`r","`n")
$ws.Range("B22").Value = $bibliografia
$ws.Range("C22").Value = $bibliografia

# --- 10. Column A should only span itself (not bleed into column B) ---
$ws.Columns.Item(1).ColumnWidth = 30.7109375
